$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns at F:G, shifting existing F,G,H,I one/two slots right (-> H,I,J,K)
$ws.Columns("F:G").Insert()

# New header labels for the two inserted columns ("Group" / "Subgroup")
$ws.Range("F7").Value = "Group"
$ws.Range("G7").Value = "Subgroup"

# Give the two new columns their own custom widths (closest representable values;
# column H keeps the width it inherited automatically from the old column F)
$ws.Columns("F").ColumnWidth = 17.666666666666668
$ws.Columns("G").ColumnWidth = 27.666666666666668

# Move the active cell selection to G8
$ws.Range("G8").Select()
